# Adjust sale/stock figures across several item rows and their
# corresponding "Sub Total" / grand-total rows in the CryCompanywiseStockReport.
# Updates quantity (F), value (G), and a handful of total cells (B), plus a
# swapped pair of rows where two batches' figures were transposed (375/376,
# 152/153).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F29").Value = 22
$ws.Range("G29").Value = 1540.88
$ws.Range("F30").Value = 86
$ws.Range("G30").Value = 3131.26
$ws.Range("F35").Value = 148
$ws.Range("G35").Value = 3742.92
$ws.Range("F42").Value = 170
$ws.Range("G42").Value = 10030
$ws.Range("F49").Value = 123
$ws.Range("G49").Value = 32069.79
$ws.Range("B54").Value = 124360.08
$ws.Range("F84").Value = 38
$ws.Range("G84").Value = 3610.76
$ws.Range("F108").Value = 21
$ws.Range("G108").Value = 2151.66
$ws.Range("B116").Value = 171961.49
$ws.Range("B152").Value = 53925
$ws.Range("B153").Value = 57756
$ws.Range("F212").Value = 2608
$ws.Range("G212").Value = 48248
$ws.Range("B217").Value = 84592.62
$ws.Range("F236").Value = 114
$ws.Range("G236").Value = 3562.5
$ws.Range("F248").Value = 0
$ws.Range("G248").Value = 0
$ws.Range("B254").Value = 98855.39
$ws.Range("F304").Value = 7
$ws.Range("G304").Value = 1398.39
$ws.Range("B313").Value = 189760.88
$ws.Range("F330").Value = 217
$ws.Range("G330").Value = 6996.08
$ws.Range("B336").Value = 35344.38
$ws.Range("F351").Value = 496
$ws.Range("G351").Value = 69732.64
$ws.Range("B353").Value = 89003.05
$ws.Range("B375").Value = 47097
$ws.Range("D375").Value = 112.28
$ws.Range("E375").Value = 134.16
$ws.Range("F375").Value = 15
$ws.Range("G375").Value = 1684.2
$ws.Range("B376").Value = 58047
$ws.Range("D376").Value = 105.54
$ws.Range("E376").Value = 126.1
$ws.Range("F376").Value = 56
$ws.Range("G376").Value = 5910.24
$ws.Range("F418").Value = 410
$ws.Range("G418").Value = 5391.5
$ws.Range("F426").Value = 375
$ws.Range("G426").Value = 6082.5
$ws.Range("B435").Value = 82084.98
$ws.Range("F459").Value = 0
$ws.Range("G459").Value = 0
$ws.Range("B463").Value = 27642
$ws.Range("F511").Value = 41
$ws.Range("G511").Value = 1339.06
$ws.Range("F515").Value = 28
$ws.Range("G515").Value = 764.96
$ws.Range("F516").Value = 35
$ws.Range("G516").Value = 6708.8
$ws.Range("F525").Value = 31
$ws.Range("G525").Value = 5347.81
$ws.Range("B526").Value = 73750.97
$ws.Range("F554").Value = 90
$ws.Range("G554").Value = 14443.2
$ws.Range("F558").Value = 25
$ws.Range("G558").Value = 2563.5
$ws.Range("F560").Value = 17
$ws.Range("G560").Value = 853.91
$ws.Range("F568").Value = 87
$ws.Range("G568").Value = 6398.85
$ws.Range("F570").Value = 161
$ws.Range("G570").Value = 11242.63
$ws.Range("F571").Value = 42
$ws.Range("G571").Value = 5953.5
$ws.Range("F575").Value = 7
$ws.Range("G575").Value = 581.14
$ws.Range("B576").Value = 160919.62
$ws.Range("F581").Value = 60
$ws.Range("G581").Value = 5697
$ws.Range("B584").Value = 41866.13
$ws.Range("F598").Value = 127
$ws.Range("G598").Value = 5483.86
$ws.Range("B599").Value = 14034.75
$ws.Range("F646").Value = 77
$ws.Range("G646").Value = 6280.12
$ws.Range("F647").Value = 110
$ws.Range("G647").Value = 5264.6
$ws.Range("F657").Value = 163
$ws.Range("G657").Value = 22006.63
$ws.Range("F658").Value = 270
$ws.Range("G658").Value = 32591.7
$ws.Range("F659").Value = 22
$ws.Range("G659").Value = 2655.62
$ws.Range("B660").Value = 117379.04
$ws.Range("F664").Value = 68
$ws.Range("G664").Value = 11144.52
$ws.Range("F672").Value = 63
$ws.Range("G672").Value = 3324.51
$ws.Range("F684").Value = 84
$ws.Range("G684").Value = 3611.16
$ws.Range("F685").Value = 149
$ws.Range("G685").Value = 7433.61
$ws.Range("F686").Value = 79
$ws.Range("G686").Value = 6367.4
$ws.Range("B690").Value = 88869.60000000001
$ws.Range("F728").Value = 2409
$ws.Range("G728").Value = 392931.99
$ws.Range("F729").Value = 287
$ws.Range("G729").Value = 81183.69
$ws.Range("F730").Value = 382
$ws.Range("G730").Value = 55256.3
$ws.Range("F734").Value = 138
$ws.Range("G734").Value = 9315
$ws.Range("B736").Value = 573537.9300000001
$ws.Range("B741").Value = 3366926.92
$ws.Range("B742").Value = 3366926.92
